$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Cells.Item(4, 8).Value = 119
$ws.Cells.Item(4, 10).Value = 99
$ws.Cells.Item(4, 12).Value = 99
$ws.Cells.Item(4, 14).Value = -327

# Row 107
$ws.Cells.Item(107, 8).Value = 855.09375
$ws.Cells.Item(107, 9).Value = 881.3929000000001
$ws.Cells.Item(107, 10).Value = 671
$ws.Cells.Item(107, 11).Value = 881.3929000000001
$ws.Cells.Item(107, 12).Value = 671
$ws.Cells.Item(107, 13).Value = 1038.6071
$ws.Cells.Item(107, 14).Value = -4511

# Row 112
$ws.Cells.Item(112, 8).Value = 2372.087
$ws.Cells.Item(112, 10).Value = 2401.388
$ws.Cells.Item(112, 12).Value = 7204.164
$ws.Cells.Item(112, 14).Value = -9420.164000000001

# Row 137
$ws.Cells.Item(137, 8).Value = 19611398
$ws.Cells.Item(137, 9).Value = 2000
$ws.Cells.Item(137, 10).Value = 25645060
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 76935180
$ws.Cells.Item(137, 13).Value = -3450
$ws.Cells.Item(137, 14).Value = -76940280

# Row 138
$ws.Cells.Item(138, 8).Value = 5330.528
$ws.Cells.Item(138, 10).Value = 6957.49
$ws.Cells.Item(138, 12).Value = 20872.47
$ws.Cells.Item(138, 14).Value = -31152.47

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 8066.887
$ws.Cells.Item(32, 9).Value = 4319.4746
$ws.Cells.Item(32, 11).Value = 4319.4746
$ws.Cells.Item(32, 13).Value = -4032.4746

# Row 61
$ws.Cells.Item(61, 8).Value = 4657.5107
$ws.Cells.Item(61, 9).Value = 4351.488
$ws.Cells.Item(61, 11).Value = 4351.488
$ws.Cells.Item(61, 13).Value = -4139.488

# Row 88
$ws.Cells.Item(88, 8).Value = 54033.934
$ws.Cells.Item(88, 9).Value = 3296
$ws.Cells.Item(88, 10).Value = 66718.414
$ws.Cells.Item(88, 11).Value = 3296
$ws.Cells.Item(88, 12).Value = 66718.414
$ws.Cells.Item(88, 13).Value = -2890
$ws.Cells.Item(88, 14).Value = -67530.414

# Row 91
$ws.Cells.Item(91, 8).Value = 54033.934
$ws.Cells.Item(91, 9).Value = 3296
$ws.Cells.Item(91, 10).Value = 66718.414
$ws.Cells.Item(91, 11).Value = 3296
$ws.Cells.Item(91, 12).Value = 66718.414
$ws.Cells.Item(91, 13).Value = -1892
$ws.Cells.Item(91, 14).Value = -69526.414

# Row 132
$ws.Cells.Item(132, 8).Value = 5202.5107
$ws.Cells.Item(132, 9).Value = 2854.3704
$ws.Cells.Item(132, 10).Value = 8372.5
$ws.Cells.Item(132, 11).Value = 8563.111199999999
$ws.Cells.Item(132, 12).Value = 25117.5
$ws.Cells.Item(132, 13).Value = -6033.111199999999
$ws.Cells.Item(132, 14).Value = -30177.5

# Row 136
$ws.Cells.Item(136, 8).Value = 4657.5107
$ws.Cells.Item(136, 9).Value = 4351.488
$ws.Cells.Item(136, 11).Value = 13054.464
$ws.Cells.Item(136, 13).Value = -10504.464

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 7897.522
$ws.Cells.Item(20, 9).Value = 7744.3076
$ws.Cells.Item(20, 11).Value = 7744.3076
$ws.Cells.Item(20, 13).Value = -7497.3076

# Row 64
$ws.Cells.Item(64, 8).Value = 611.1905
$ws.Cells.Item(64, 10).Value = 820.8182
$ws.Cells.Item(64, 12).Value = 820.8182
$ws.Cells.Item(64, 14).Value = -1270.8182

# Row 67
$ws.Cells.Item(67, 8).Value = 611.1905
$ws.Cells.Item(67, 10).Value = 820.8182
$ws.Cells.Item(67, 12).Value = 820.8182
$ws.Cells.Item(67, 14).Value = -2380.8182

# Row 128
$ws.Cells.Item(128, 8).Value = 7608.75
$ws.Cells.Item(128, 9).Value = 7608.75
$ws.Cells.Item(128, 11).Value = 22826.25
$ws.Cells.Item(128, 13).Value = -20336.25

# Row 140
$ws.Cells.Item(140, 8).Value = 188968.42
$ws.Cells.Item(140, 10).Value = 188968.42
$ws.Cells.Item(140, 12).Value = 188968.42
$ws.Cells.Item(140, 14).Value = -199328.42

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 124.71429
$ws.Cells.Item(7, 9).Value = 79.8
$ws.Cells.Item(7, 11).Value = 79.8
$ws.Cells.Item(7, 13).Value = 33.2

# Row 16
$ws.Cells.Item(16, 8).Value = 1187
$ws.Cells.Item(16, 9).Value = 999.2
$ws.Cells.Item(16, 11).Value = 999.2
$ws.Cells.Item(16, 13).Value = -712.2

# Row 31
$ws.Cells.Item(31, 8).Value = 26318940
$ws.Cells.Item(31, 9).Value = 35716230
$ws.Cells.Item(31, 10).Value = 6521.4
$ws.Cells.Item(31, 11).Value = 35716230
$ws.Cells.Item(31, 12).Value = 6521.4
$ws.Cells.Item(31, 13).Value = -35715935
$ws.Cells.Item(31, 14).Value = -7111.4

# Row 34
$ws.Cells.Item(34, 8).Value = 26318940
$ws.Cells.Item(34, 9).Value = 35716230
$ws.Cells.Item(34, 10).Value = 6521.4
$ws.Cells.Item(34, 11).Value = 35716230
$ws.Cells.Item(34, 12).Value = 6521.4
$ws.Cells.Item(34, 13).Value = -35716028
$ws.Cells.Item(34, 14).Value = -6925.4

# Row 62
$ws.Cells.Item(62, 8).Value = 35883.363
$ws.Cells.Item(62, 9).Value = 24535.666
$ws.Cells.Item(62, 10).Value = 60199.855
$ws.Cells.Item(62, 11).Value = 24535.666
$ws.Cells.Item(62, 12).Value = 60199.855
$ws.Cells.Item(62, 13).Value = -23911.666
$ws.Cells.Item(62, 14).Value = -61447.855

# Row 65
$ws.Cells.Item(65, 8).Value = 35883.363
$ws.Cells.Item(65, 9).Value = 24535.666
$ws.Cells.Item(65, 10).Value = 60199.855
$ws.Cells.Item(65, 11).Value = 122678.33
$ws.Cells.Item(65, 12).Value = 300999.275
$ws.Cells.Item(65, 13).Value = -119558.33
$ws.Cells.Item(65, 14).Value = -307239.275

# Row 97
$ws.Cells.Item(97, 8).Value = 10293.941
$ws.Cells.Item(97, 10).Value = 10749.8125
$ws.Cells.Item(97, 12).Value = 10749.8125
$ws.Cells.Item(97, 14).Value = -12731.8125

# Row 105
$ws.Cells.Item(105, 8).Value = 928
$ws.Cells.Item(105, 9).Value = 611.4286
$ws.Cells.Item(105, 11).Value = 611.4286
$ws.Cells.Item(105, 13).Value = 1135.5714

# Row 113
$ws.Cells.Item(113, 8).Value = 1187
$ws.Cells.Item(113, 9).Value = 999.2
$ws.Cells.Item(113, 11).Value = 999.2
$ws.Cells.Item(113, 13).Value = 1170.8

# Row 122
$ws.Cells.Item(122, 8).Value = 3983.2307
$ws.Cells.Item(122, 9).Value = 2616.7144
$ws.Cells.Item(122, 11).Value = 7850.1432
$ws.Cells.Item(122, 13).Value = -5400.1432

# Row 132
$ws.Cells.Item(132, 8).Value = 58826548
$ws.Cells.Item(132, 9).Value = 66668956
$ws.Cells.Item(132, 10).Value = 8494.5
$ws.Cells.Item(132, 11).Value = 200006868
$ws.Cells.Item(132, 12).Value = 25483.5
$ws.Cells.Item(132, 13).Value = -200004338
$ws.Cells.Item(132, 14).Value = -30543.5

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Cells.Item(33, 8).Value = 457.77777
$ws.Cells.Item(33, 10).Value = 141.66667
$ws.Cells.Item(33, 12).Value = 850.0000200000001
$ws.Cells.Item(33, 14).Value = -1416.00002

# Row 37
$ws.Cells.Item(37, 8).Value = 100101976
$ws.Cells.Item(37, 10).Value = 100101976
$ws.Cells.Item(37, 12).Value = 300305928
$ws.Cells.Item(37, 14).Value = -300306152

# Row 131
$ws.Cells.Item(131, 8).Value = 19940034
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 13).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 5321.8
$ws.Cells.Item(134, 9).Value = 4802
$ws.Cells.Item(134, 10).Value = 10000
$ws.Cells.Item(134, 11).Value = 14406
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 13).Value = -9336
$ws.Cells.Item(134, 14).Value = -40140

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 822
$ws.Cells.Item(107, 9).Value = 437
$ws.Cells.Item(107, 10).Value = 1688.25
$ws.Cells.Item(107, 11).Value = 437
$ws.Cells.Item(107, 12).Value = 1688.25
$ws.Cells.Item(107, 13).Value = 1483
$ws.Cells.Item(107, 14).Value = -5528.25

# Row 122
$ws.Cells.Item(122, 8).Value = 4408.2334
$ws.Cells.Item(122, 9).Value = 1906.7142
$ws.Cells.Item(122, 11).Value = 5720.142599999999
$ws.Cells.Item(122, 13).Value = -3270.142599999999

# Row 123
$ws.Cells.Item(123, 8).Value = 51587.637
$ws.Cells.Item(123, 10).Value = 51587.637
$ws.Cells.Item(123, 12).Value = 51587.637
$ws.Cells.Item(123, 14).Value = -56487.637

# Row 141
$ws.Cells.Item(141, 8).Value = 84982.336
$ws.Cells.Item(141, 10).Value = 84982.336
$ws.Cells.Item(141, 12).Value = 84982.336
$ws.Cells.Item(141, 14).Value = -95342.336

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5048.625
$ws.Cells.Item(7, 10).Value = 5565
$ws.Cells.Item(7, 12).Value = 5565
$ws.Cells.Item(7, 14).Value = -5789

# Row 82
$ws.Cells.Item(82, 8).Value = 4810.8823
$ws.Cells.Item(82, 9).Value = 3746
$ws.Cells.Item(82, 10).Value = 5556.3
$ws.Cells.Item(82, 11).Value = 3746
$ws.Cells.Item(82, 12).Value = 5556.3
$ws.Cells.Item(82, 13).Value = -3385
$ws.Cells.Item(82, 14).Value = -6278.3

# Row 85
$ws.Cells.Item(85, 8).Value = 4810.8823
$ws.Cells.Item(85, 9).Value = 3746
$ws.Cells.Item(85, 10).Value = 5556.3
$ws.Cells.Item(85, 11).Value = 3746
$ws.Cells.Item(85, 12).Value = 5556.3
$ws.Cells.Item(85, 13).Value = -2498
$ws.Cells.Item(85, 14).Value = -8052.3

# Row 126
$ws.Cells.Item(126, 8).Value = 5048.625
$ws.Cells.Item(126, 10).Value = 5565
$ws.Cells.Item(126, 12).Value = 16695
$ws.Cells.Item(126, 14).Value = -21635

# Row 136
$ws.Cells.Item(136, 8).Value = 5611.75
$ws.Cells.Item(136, 9).Value = 3535.25
$ws.Cells.Item(136, 10).Value = 6650
$ws.Cells.Item(136, 11).Value = 10605.75
$ws.Cells.Item(136, 12).Value = 19950
$ws.Cells.Item(136, 13).Value = -8055.75
$ws.Cells.Item(136, 14).Value = -25050

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Cells.Item(14, 8).Value = 4737.5

# Row 74
$ws.Cells.Item(74, 8).Value = 9918.909
$ws.Cells.Item(74, 10).Value = 9918.909
$ws.Cells.Item(74, 12).Value = 9918.909
$ws.Cells.Item(74, 14).Value = -11790.909

# Row 77
$ws.Cells.Item(77, 8).Value = 9918.909
$ws.Cells.Item(77, 10).Value = 9918.909
$ws.Cells.Item(77, 12).Value = 29756.727
$ws.Cells.Item(77, 14).Value = -39116.727

# Row 81
$ws.Cells.Item(81, 8).Value = 7597.846
$ws.Cells.Item(81, 9).Value = 5510.3335
$ws.Cells.Item(81, 10).Value = 10444.454
$ws.Cells.Item(81, 11).Value = 11020.667
$ws.Cells.Item(81, 12).Value = 20888.908
$ws.Cells.Item(81, 13).Value = -9959.666999999999
$ws.Cells.Item(81, 14).Value = -23010.908

# Row 84
$ws.Cells.Item(84, 8).Value = 7597.846
$ws.Cells.Item(84, 9).Value = 5510.3335
$ws.Cells.Item(84, 10).Value = 10444.454
$ws.Cells.Item(84, 11).Value = 55103.335
$ws.Cells.Item(84, 12).Value = 104444.54
$ws.Cells.Item(84, 13).Value = -49799.335
$ws.Cells.Item(84, 14).Value = -115052.54

# Row 126
$ws.Cells.Item(126, 8).Value = 500000500
$ws.Cells.Item(126, 9).Value = 500000500
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 1500001500
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -1499999030
$ws.Cells.Item(126, 14).ClearContents()

# Row 132
$ws.Cells.Item(132, 8).Value = 4245
$ws.Cells.Item(132, 9).Value = 1916.2
$ws.Cells.Item(132, 11).Value = 5748.6
$ws.Cells.Item(132, 13).Value = -3218.6

# Row 136
$ws.Cells.Item(136, 8).Value = 4469.95
$ws.Cells.Item(136, 9).Value = 2558.3809
$ws.Cells.Item(136, 11).Value = 7675.1427
$ws.Cells.Item(136, 13).Value = -5125.1427
